$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.007.77"
$ws.Range("E2").Value = "  +3.19%  "
$ws.Range("D3").Value = "3.284.75"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "630.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.384"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +20.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.690"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +16.01%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "3.284.84"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.577"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.48%  "
$ws.Range("E12").Value = "  +12.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "3.887.37"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "88.937.24"
$ws.Range("E17").Value = "  +3.68%  "
$ws.Range("D18").Value = "3.290.20"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "438.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("E23").Value = "  +2.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("B25").Value = "Aptos"
$ws.Range("C25").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.92%  "
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "76.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000136"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.185"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.85%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "567.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.44%  "
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("E38").Value = "  -7.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("E43").Value = "  -3.58%  "
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "155.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "181.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("E50").Value = "  +14.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0679"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +21.83%  "
